$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date column (A2:A6) from 2025-12-05 to 2025-12-06.
# Force text format first so Excel doesn't auto-convert the string into a
# date serial number, then restore the default (Normal) style so the cell
# keeps its original (unstyled) appearance.
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("A2:A6").Value = "2025-12-06"
$ws.Range("A2:A6").Style = "Normal"

# Row 4 and Row 5 swap identity: row4 becomes MARA Holdings / MARA,
# row5 becomes Coinbase Global, Inc. / COIN
$ws.Range("B4").Value = "MARA Holdings, Inc."
$ws.Range("C4").Value = "MARA"
$ws.Range("B5").Value = "Coinbase Global, Inc."
$ws.Range("C5").Value = "COIN"

# Row 2 (RIOT) updated metrics
$ws.Range("D2").Value = 15.27
$ws.Range("E2").Value = 60.2
$ws.Range("F2").Value = -5.33
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 53
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 54.5
$ws.Range("N2").Value = 51.54219175917372

# Row 3 (BTC-USD) updated metrics
$ws.Range("D3").Value = 91134.45
$ws.Range("E3").Value = 62.9
$ws.Range("F3").Value = 0.82
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 53
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 53.5
$ws.Range("N3").Value = 51.54219175917372

# Row 4 (now MARA) updated metrics
$ws.Range("D4").Value = 12.09
$ws.Range("E4").Value = 50.8
$ws.Range("F4").Value = 2.33
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 63
$ws.Range("J4").Value = 73
$ws.Range("K4").Value = 52.7
$ws.Range("N4").Value = 51.54219175917372

# Row 5 (now COIN) updated metrics
$ws.Range("D5").Value = 270.9
$ws.Range("E5").Value = 44.6
$ws.Range("F5").Value = -0.7
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 56
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 48.5
$ws.Range("N5").Value = 51.54219175917372

# Row 6 (MSTR) updated metrics
$ws.Range("D6").Value = 183.17
$ws.Range("E6").Value = 41.7
$ws.Range("F6").Value = 3.38
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 36
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 36
$ws.Range("K6").Value = 43.5
$ws.Range("N6").Value = 51.54219175917372
